# Apply data update to the "Inscricoes" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 5: Inscritos (E5) 31 -> 32
$ws.Range("E5").Value = 32

# Row 6: Inscritos (E6) 58 -> 59, Pagos (F6) 28 -> 29, Inscricoes homologadas (H6) 35 -> 36
$ws.Range("E6").Value = 59
$ws.Range("F6").Value = 29
$ws.Range("H6").Value = 36

# Row 8: Inscritos (E8) 49 -> 50
$ws.Range("E8").Value = 50

# Row 16: Inscritos (E16) 330 -> 332, Pagos (F16) 95 -> 96, Inscricoes homologadas (H16) 183 -> 184
$ws.Range("E16").Value = 332
$ws.Range("F16").Value = 96
$ws.Range("H16").Value = 184
